# The workbook was re-exported and the "Förändrad" (Changed) date in
# column C was bumped by one day (45203 -> 45204) for every data row
# (rows 2 through 74) on the single worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("C2:C74")
foreach ($cell in $rng.Cells) {
    if ($cell.Value2 -eq 45203) {
        $cell.Value2 = 45204
    }
}
